$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price values must stay as text (matches source "t=inlineStr" cells),
# so force Text format before assigning, then restore Normal style (no explicit s=) after.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4067"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08221"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.064"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.247"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06893"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.189"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.641"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.117"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.017"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09613"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.655"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.545"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06103"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.184"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.066"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5973"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.386"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07599"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5593"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.426"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.31"
$ws.Range("D51").Style = "Normal"

# Remaining text cells (names, links, multi-dot prices, percentages)
$ws.Range("D2").Value = "29.437.21"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.918.34"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.930.47"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "29.457.02"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "2.166.15"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +9.06%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +6.37%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("E51").Value = "  +0.06%  "
